$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Run ID value (B2) with the corrected translation GUID
$ws.Range("B2").Value = "5f85f5bb-596e-44e1-b790-414f94ed8f18_5"

# Updated accuracy-ish metrics for rows 3-23 (column B)
$ws.Range("B3").Value = 0.68571000000000004
$ws.Range("B4").Value = 0.65832999999999997
$ws.Range("B5").Value = 0.70204
$ws.Range("B6").Value = 0.65832999999999997
$ws.Range("B7").Value = 0.63402999999999998
$ws.Range("B8").Value = 0.71223999999999998
$ws.Range("B9").Value = 0.73980000000000001
$ws.Range("B10").Value = 0.65832999999999997
$ws.Range("B11").Value = 0.56384000000000001
$ws.Range("B12").Value = 0.68571000000000004
$ws.Range("B13").Value = 0.64136000000000004
$ws.Range("B14").Value = 0.64990999999999999
$ws.Range("B15").Value = 0.31667000000000001
$ws.Range("B16").Value = 0.31667000000000001
$ws.Range("B17").Value = 0.60119
$ws.Range("B18").Value = 0.68571000000000004
$ws.Range("B19").Value = 0.73129
$ws.Range("B20").Value = 0.65832999999999997
$ws.Range("B21").Value = 0.68571000000000004
$ws.Range("B22").Value = 0.68571000000000004
$ws.Range("B23").Value = 0.70443

# Update the current selection to reflect the new active range
$ws.Range("A3:A23").Select()
